$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for RM Part No L/R/O
$ws.Range("D1").Value = "RM Part No L"
$ws.Range("E1").Value = "RM Part No R"
$ws.Range("F1").Value = "RM Part No O"

# E1: same bold style as A1:C1 (fontId 1), just add vertical centering
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E1").VerticalAlignment = -4108

# D1 and F1: new bold 12pt font (not inheriting theme minor scheme), vertical centered
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Size = 12
$ws.Range("D1").VerticalAlignment = -4108

$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Font.Size = 12
$ws.Range("F1").VerticalAlignment = -4108

# Column widths
$ws.Columns.Item(4).ColumnWidth = 21.33203125
$ws.Columns.Item(5).ColumnWidth = 18.77734375
$ws.Columns.Item(6).ColumnWidth = 19.21875

# Selection on F1 to match final state
$ws.Range("F1").Select()
